# "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to model a single "Property" was renamed to
# "DataNode" to line up with the rest of the data-config tables
# (DataTable / Entity) that share the same naming scheme. Re-apply the
# same interactive edits a user made in Excel: rename the tab, leave the
# cursor where they left it (B41), and nudge the two text columns a hair
# wider (a side effect of the same editing session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Re-select the cell the editor ended up on (was K17, now B41).
$ws.Range("B41").Select() | Out-Null

# Small column width tweaks for column A ("Id"/name column) and column C
# (third data column) made during the same pass.
$ws.Columns.Item(1).ColumnWidth = 31.142857142857142
$ws.Columns.Item(3).ColumnWidth = 30.857142857142858
